$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the value to be stored as text, preserving exact formatting
    # (e.g. trailing zeros / dotted thousand separators) instead of letting
    # Excel auto-convert number-looking strings into numeric cells.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "40.976.10"
$ws.Range("E2").Value = "  -6.77%  "
Set-TextValue $ws.Range("D3") "2.182.08"
$ws.Range("E3").Value = "  -7.49%  "
$ws.Range("E4").Value = "  -0.19%  "
Set-TextValue $ws.Range("D5") "239.59"
$ws.Range("E5").Value = "  -0.24%  "
Set-TextValue $ws.Range("D6") "0.618"
$ws.Range("E6").Value = "  -7.47%  "
Set-TextValue $ws.Range("D7") "69.63"
$ws.Range("E7").Value = "  -5.09%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -12.35%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D10") "36.22"
$ws.Range("E10").Value = "  +6.81%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D11") "57.61"
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("E12").Value = "  -8.46%  "
$ws.Range("E13").Value = "  -4.35%  "
$ws.Range("E14").Value = "  -9.64%  "
Set-TextValue $ws.Range("D15") "2.503.46"
$ws.Range("E15").Value = "  -7.63%  "
Set-TextValue $ws.Range("D16") "14.50"
$ws.Range("E16").Value = "  -10.48%  "
Set-TextValue $ws.Range("D17") "0.827"
$ws.Range("E17").Value = "  -9.05%  "
Set-TextValue $ws.Range("D18") "2.181.46"
$ws.Range("E18").Value = "  -7.55%  "
Set-TextValue $ws.Range("D19") "40.986.67"
$ws.Range("E19").Value = "  -6.69%  "
Set-TextValue $ws.Range("D20") "0.0₃0932"
$ws.Range("E20").Value = "  -9.33%  "
Set-TextValue $ws.Range("D21") "73.16"
$ws.Range("E21").Value = "  -5.91%  "
Set-TextValue $ws.Range("D22") "5.99"
$ws.Range("E22").Value = "  -8.47%  "
Set-TextValue $ws.Range("D23") "229.85"
$ws.Range("E23").Value = "  -9.15%  "
$ws.Range("E24").Value = "  +7.27%  "
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("E28").Value = "  -5.09%  "
Set-TextValue $ws.Range("D29") "9.63"
$ws.Range("E29").Value = "  -7.94%  "
Set-TextValue $ws.Range("D30") "166.61"
$ws.Range("E30").Value = "  -5.55%  "
Set-TextValue $ws.Range("D31") "20.08"
$ws.Range("E31").Value = "  -9.84%  "
$ws.Range("E32").Value = "  -9.58%  "
Set-TextValue $ws.Range("D33") "0.123"
$ws.Range("E33").Value = "  -8.24%  "
Set-TextValue $ws.Range("D34") "0.0695"
$ws.Range("E34").Value = "  -6.99%  "
Set-TextValue $ws.Range("D35") "5.05"
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("E36").Value = "  -10.62%  "
$ws.Range("E37").Value = "  +0.90%  "
Set-TextValue $ws.Range("D38") "23.34"
$ws.Range("E38").Value = "  +17.88%  "
Set-TextValue $ws.Range("D39") "2.23"
$ws.Range("E39").Value = "  -7.38%  "
$ws.Range("E40").Value = "  -3.48%  "
Set-TextValue $ws.Range("D41") "5.73"
$ws.Range("E41").Value = "  -13.16%  "
Set-TextValue $ws.Range("D42") "64.73"
$ws.Range("E42").Value = "  -0.26%  "
Set-TextValue $ws.Range("D43") "4.79"
$ws.Range("E43").Value = "  -12.03%  "
Set-TextValue $ws.Range("D44") "8.59"
$ws.Range("E44").Value = "  -5.10%  "
Set-TextValue $ws.Range("D45") "0.189"
$ws.Range("E45").Value = "  -5.96%  "
Set-TextValue $ws.Range("D47") "0.0977"
$ws.Range("E47").Value = "  -8.39%  "
Set-TextValue $ws.Range("D48") "4.49"
$ws.Range("E48").Value = "  +4.36%  "
Set-TextValue $ws.Range("D49") "9.96"
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("E50").Value = "  -6.61%  "
Set-TextValue $ws.Range("D51") "1.08"
$ws.Range("E51").Value = "  -6.77%  "
